$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.613.39'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.739.44'
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.36'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.0000'
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4938'
$ws.Range("E7").Value = '  +2.86%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2675'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06278'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.744.63'
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07050'
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6143'
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.584'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '78.04'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.0000'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.635.97'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9999'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007273'
$ws.Range("E19").Value = '  +4.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.57'
$ws.Range("E20").Value = '  -1.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.971.08'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.564'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.719'
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.275'
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.88'
$ws.Range("E25").Value = '  +1.78%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '107.44'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.754'
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08044'
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.732'
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04623'
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.610'
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.014'
$ws.Range("E36").Value = '  +2.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6396'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.064'
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9036'
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.427'
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.95'
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.440'
$ws.Range("E44").Value = '  -4.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3934'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.867'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1185'
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05390'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.63'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.786'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("E51").Value = '  -0.99%  '
